# Rename the worksheet currently named "Sheet1" to "Research".
#
# The workbook's tab order is:
#   1 About
#   2 Metadata
#   3 Information Model
#   4 Sheet1                 <- rename this one to "Research"
#   5 Data
#   6 SNAQ65+AppetiteScoreCodelist
#   7 SNAQ65+ExerciseScoreCodelist
#   8 SNAQ65+UpperarmCircumferenceCod
#   9 SNAQ65+WeightLossScoreCodelist
#  10 Terms of Use

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Research"
